# feat: add 2022-Q4 data
#
# - Insert a new "2022-Q4" sheet (with its fund holdings data) positioned
#   right after "总计" and before the existing "2022-Q2" sheet.
# - Insert a new summary row on "总计" for 2022-Q4, pushing the existing
#   2022-Q2 summary row down to row 3 (and bumping its index column).

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)     # "总计"
$q2      = $wb.Worksheets.Item(2)     # existing "2022-Q2"

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: push the existing 2022-Q2 row
#    down to row 3, and write the new 2022-Q4 row into row 2.
# ---------------------------------------------------------------------
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.09

# Copy the index-column formatting (bold/border style) from A2 onto A3.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.14

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet before the existing "2022-Q2"
#    sheet (so tab order becomes 总计, 2022-Q4, 2022-Q2).
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($q2)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2 - fund 014277. Columns B, D, E, F, G look numeric, so force text
# formatting before assigning them or Excel will strip the leading zero
# / coerce the string to a number.
$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "014277"
$q4.Range("C2").Value = "万家北交所慧选两年定期开放混合A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "3.25"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "94.43"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "3.65"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.1186"
$q4.Range("H2").Value = 10

# Row 3 - fund 014278.
$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "014278"
$q4.Range("C3").Value = "万家北交所慧选两年定期开放混合C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.45"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "94.43"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "3.65"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0164"
$q4.Range("H3").Value = 10

# Drop the temporary "@" text format again so the data cells end up on
# the plain/default style, then (re)apply the shared bold+border style
# (the one already used by the "总计" sheet's header row / index column)
# to the new sheet's header row and index column only.
$q4.Range("B2:B3").ClearFormats()
$q4.Range("D2:G3").ClearFormats()

$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)     # xlPasteFormats

$summary.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)     # xlPasteFormats

# Keep the original "2022-Q2" sheet as the active/selected tab (it was the
# active sheet before this edit; only its tab position moved). Re-resolve
# it by name since $q2 now points at a stale position after the insert.
$wb.Worksheets.Item("2022-Q2").Activate()

Write-Output "2022-Q4 sheet added"
